# update and run DE density calcs
# Adds two new "basic" summary sheets (area_lores_basic, area_pop_sum_basic)
# mirroring the existing area_lores / area_pop_sum sheets but with refreshed
# values, placed after the existing sheets.

$wb = $excel.ActiveWorkbook

$srcLores = $wb.Worksheets.Item("area_lores")
$srcPopSum = $wb.Worksheets.Item("area_pop_sum")

# ----------------------------------------------------------------------
# new sheet: area_lores_basic
# ----------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsLoresBasic = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsLoresBasic.Name = "area_lores_basic"

# match sheetPr/outlinePr + pageMargins used on the original sheets
$wsLoresBasic.Outline.SummaryRow = 1        # xlSummaryBelow
$wsLoresBasic.Outline.SummaryColumn = -4152 # xlSummaryOnRight
$wsLoresBasic.PageSetup.LeftMargin = 0.75 * 72
$wsLoresBasic.PageSetup.RightMargin = 0.75 * 72
$wsLoresBasic.PageSetup.TopMargin = 1 * 72
$wsLoresBasic.PageSetup.BottomMargin = 1 * 72
$wsLoresBasic.PageSetup.HeaderMargin = 0.5 * 72
$wsLoresBasic.PageSetup.FooterMargin = 0.5 * 72

$wsLoresBasic.Range("A1").Value = "index"
$wsLoresBasic.Range("B1").Value = "area"
$wsLoresBasic.Range("A2").Value = "count"
$wsLoresBasic.Range("B2").Value = 34
$wsLoresBasic.Range("A3").Value = "mean"
$wsLoresBasic.Range("B3").Value = 8.761330811664413
$wsLoresBasic.Range("A4").Value = "std"
$wsLoresBasic.Range("B4").Value = 8.469738065308816
$wsLoresBasic.Range("A5").Value = "min"
$wsLoresBasic.Range("B5").Value = 1.643797814826138
# these three labels look like percentages ("25%", "50%", "75%") -- force
# literal text entry (quote-prefix) so Excel doesn't reinterpret them as
# numeric percent values
$wsLoresBasic.Range("A6").Formula = "'25%"
$wsLoresBasic.Range("B6").Value = 3.219074674173791
$wsLoresBasic.Range("A7").Formula = "'50%"
$wsLoresBasic.Range("B7").Value = 5.487450795955053
$wsLoresBasic.Range("A8").Formula = "'75%"
$wsLoresBasic.Range("B8").Value = 9.216428102088065
$wsLoresBasic.Range("A9").Value = "max"
$wsLoresBasic.Range("B9").Value = 37.96161634642002

# strip the quote-prefix/text formatting picked up above back to plain
# general formatting (matches the plain, unstyled data cells elsewhere)
$wsLoresBasic.Range("A5").Copy()
$wsLoresBasic.Range("A6:A8").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# match header styling (bold, centered, bordered) used on area_lores!A1:B1
$srcLores.Range("A1:B1").Copy()
$wsLoresBasic.Range("A1:B1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# ----------------------------------------------------------------------
# new sheet: area_pop_sum_basic
# ----------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPopSumBasic = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsPopSumBasic.Name = "area_pop_sum_basic"

# match sheetPr/outlinePr + pageMargins used on the original sheets
$wsPopSumBasic.Outline.SummaryRow = 1        # xlSummaryBelow
$wsPopSumBasic.Outline.SummaryColumn = -4152 # xlSummaryOnRight
$wsPopSumBasic.PageSetup.LeftMargin = 0.75 * 72
$wsPopSumBasic.PageSetup.RightMargin = 0.75 * 72
$wsPopSumBasic.PageSetup.TopMargin = 1 * 72
$wsPopSumBasic.PageSetup.BottomMargin = 1 * 72
$wsPopSumBasic.PageSetup.HeaderMargin = 0.5 * 72
$wsPopSumBasic.PageSetup.FooterMargin = 0.5 * 72

$wsPopSumBasic.Range("A1").Value = "index"
$wsPopSumBasic.Range("B1").Value = 0
$wsPopSumBasic.Range("A2").Value = "area"
$wsPopSumBasic.Range("B2").Value = 297.88524759659
$wsPopSumBasic.Range("A3").Value = "population"
$wsPopSumBasic.Range("B3").Value = 503235
$wsPopSumBasic.Range("A4").Value = "density"
$wsPopSumBasic.Range("B4").Value = 1689.358583750693

# match header styling (bold, centered, bordered) used on area_pop_sum!A1:B1
$srcPopSum.Range("A1:B1").Copy()
$wsPopSumBasic.Range("A1:B1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# restore original active sheet/selection (area_lores) like before.xlsx
$srcLores.Activate()

Write-Host "Added sheets: $($wb.Worksheets.Count) total"
